# Weekly update: insert a new price-observation row for Arándano (blue) /
# Vega Modelo de Temuco at row 54, pushing the existing rows 54-110 down
# to 55-111 (dates keep their original relative order beneath the new entry).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$row = 54
$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 44897
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 430
$ws.Cells.Item($row, 14).Value = 2000
$ws.Cells.Item($row, 15).Value = 2200
$ws.Cells.Item($row, 16).Value = 2116
$ws.Cells.Item($row, 17).Value = "`$/kilo"
$ws.Cells.Item($row, 18).Value = "Región del Maule"
$ws.Cells.Item($row, 19).Value = 2116
$ws.Cells.Item($row, 20).Value = 1
